$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.852.84"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").Value = "2.276.25"
$ws.Range("E3").Value = "  -3.03%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'314.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "'101.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.20%  "
$ws.Range("D7").Value = "'0.624"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.15%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.599"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.82%  "
$ws.Range("D10").Value = "'38.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.39%  "
$ws.Range("D11").Value = "'0.0899"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.82%  "
$ws.Range("D12").Value = "'8.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.58%  "
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "'0.955"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.11%  "
$ws.Range("D15").Value = "'15.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.06%  "
$ws.Range("D16").Value = "2.620.96"
$ws.Range("E16").Value = "  -3.07%  "
$ws.Range("D17").Value = "2.282.16"
$ws.Range("E17").Value = "  -4.65%  "
$ws.Range("D18").Value = "41.802.09"
$ws.Range("E18").Value = "  -1.78%  "
$ws.Range("D19").Value = "'7.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").Value = "'282.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.57%  "
$ws.Range("D22").Value = "'73.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").Value = "'3.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.42%  "
$ws.Range("D24").Value = "'2.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("D25").Value = "'9.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.42%  "
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").Value = "'10.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.13%  "
$ws.Range("D29").Value = "'22.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D30").Value = "'162.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.02%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.0872"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.99%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "'34.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.00%  "
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").Value = "'5.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.50%  "
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").Value = "'0.115"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.26%  "
$ws.Range("D37").Value = "'4.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.13%  "
$ws.Range("D38").Value = "'2.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.05%  "
$ws.Range("D39").Value = "'0.0344"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.01%  "
$ws.Range("D40").Value = "'3.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.70%  "
$ws.Range("D41").Value = "'102.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +21.38%  "
$ws.Range("D42").Value = "'1.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("D43").Value = "'68.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").Value = "'0.223"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.09%  "
$ws.Range("D46").Value = "'114.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.61%  "
$ws.Range("D47").Value = "'11.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.80%  "
$ws.Range("D48").Value = "'8.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.64%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'5.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").Value = "'75.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").Value = "'1.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.37%  "
